# Apply the commit's edits to the "data_sources" workbook:
#  - Update the WLC note in B4 with the added copyright caveat for ESV/NRSV.
#  - Add a new note in E2 documenting the Claude-assisted CSV.
#  - Move the active selection to B5 (as last left by the editing session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Claude 4.1 for locations_no_geo.csv"

$ws.Range("B4").Value = 'wlc = "Westminster Leningrad Codex" (oldest known complete Hebrew manuscript); ESV and NRSV cannot be downloaded locally per copyright'

$ws.Range("B5").Select()
